$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.715993262238555
$ws.Range("C2").Value = 0.1559993107592845
$ws.Range("D2").Value = 0.1085831625462568
$ws.Range("F2").Value = 1.899826209119951
$ws.Range("G2").Value = 1.280560895383985
$ws.Range("H2").Value = 1.201362203969296
$ws.Range("I2").Value = 1.088450046307784
$ws.Range("J2").Value = 0.1550793913403208
$ws.Range("L2").Value = 0.4332628079233842
$ws.Range("N2").Value = 1.462311521678892
$ws.Range("B3").Value = 1.611428143844762
$ws.Range("C3").Value = 0.1363617349252024
$ws.Range("D3").Value = 0.1083094924820962
$ws.Range("F3").Value = 1.897486873322478
$ws.Range("G3").Value = 1.274352606944674
$ws.Range("H3").Value = 1.204307485961664
$ws.Range("I3").Value = 1.094421106830453
$ws.Range("J3").Value = 0.1562286021713
$ws.Range("L3").Value = 0.4256806792747767
$ws.Range("N3").Value = 1.479129206295656
$ws.Range("B4").Value = 1.547812373119939
$ws.Range("C4").Value = 0.12426083606681
$ws.Range("D4").Value = 0.1081597207920346
$ws.Range("F4").Value = 1.89714534110206
$ws.Range("G4").Value = 1.271416668153634
$ws.Range("H4").Value = 1.206735315323556
$ws.Range("I4").Value = 1.098706173385189
$ws.Range("J4").Value = 0.1569795039546129
$ws.Range("L4").Value = 0.4212082527224652
$ws.Range("N4").Value = 1.490022803541706
$ws.Range("B5").Value = 1.522037332755701
$ws.Range("C5").Value = 0.1193187522621599
$ws.Range("D5").Value = 0.1081033038809736
$ws.Range("F5").Value = 1.897281242029749
$ws.Range("G5").Value = 1.270440085964978
$ws.Range("H5").Value = 1.207880364352974
$ws.Range("I5").Value = 1.100607806551061
$ws.Range("J5").Value = 0.1572969020759825
$ws.Range("L5").Value = 0.419431850053229
$ws.Range("N5").Value = 1.494604712924168
$ws.Range("B6").Value = 1.517766434005978
$ws.Range("C6").Value = 0.1184974674899308
$ws.Range("D6").Value = 0.1080942154144822
$ws.Range("F6").Value = 1.897320415796074
$ws.Range("G6").Value = 1.270291190180757
$ws.Range("H6").Value = 1.208079900156875
$ws.Range("I6").Value = 1.100932954122882
$ws.Range("J6").Value = 0.1573502946831358
$ws.Range("L6").Value = 0.4191396705770103
$ws.Range("N6").Value = 1.495374152661533
$ws.Range("B7").Value = 1.547464157500769
$ws.Range("C7").Value = 0.1241942293616205
$ws.Range("D7").Value = 0.108158941210295
$ws.Range("F7").Value = 1.897146060432476
$ws.Range("G7").Value = 1.271402608092899
$ws.Range("H7").Value = 1.20675012757188
$ws.Range("I7").Value = 1.098731190360603
$ws.Range("J7").Value = 0.1569837383228396
$ws.Range("L7").Value = 0.4211841085209471
$ws.Range("N7").Value = 1.490084019089949
$ws.Range("B8").Value = 1.679817890517427
$ws.Range("C8").Value = 0.1492372840313578
$ws.Range("D8").Value = 0.1084850247693048
$ws.Range("F8").Value = 1.898792231587919
$ws.Range("G8").Value = 1.278238201088342
$ws.Range("H8").Value = 1.202249134077462
$ws.Range("I8").Value = 1.090380328983436
$ws.Range("J8").Value = 0.1554662489510417
$ws.Range("L8").Value = 0.4306105719666533
$ws.Range("N8").Value = 1.467992472116475
$ws.Range("B9").Value = 1.943990447234398
$ws.Range("C9").Value = 0.1980037218772281
$ws.Range("D9").Value = 0.1092684120551297
$ws.Range("F9").Value = 1.910720507070835
$ws.Range("G9").Value = 1.298615559723487
$ws.Range("H9").Value = 1.198341752959138
$ws.Range("I9").Value = 1.078922766347084
$ws.Range("J9").Value = 0.1528491416814308
$ws.Range("L9").Value = 0.4505441267586434
$ws.Range("N9").Value = 1.429173905298999
$ws.Range("B10").Value = 2.140873084108989
$ws.Range("C10").Value = 0.233628293650554
$ws.Range("D10").Value = 0.1099305341553958
$ws.Range("F10").Value = 1.924810570993145
$ws.Range("G10").Value = 1.317872398671511
$ws.Range("H10").Value = 1.198477470917794
$ws.Range("I10").Value = 1.073516088255431
$ws.Range("J10").Value = 0.1511441597795251
$ws.Range("L10").Value = 0.466068966495925
$ws.Range("N10").Value = 1.403398749473769
$ws.Range("B11").Value = 2.231043002134584
$ws.Range("C11").Value = 0.2497917851735281
$ws.Range("D11").Value = 0.1102503082797455
$ws.Range("F11").Value = 1.932382257568406
$ws.Range("G11").Value = 1.327571184001215
$ws.Range("H11").Value = 1.199193851426685
$ws.Range("I11").Value = 1.0717131593102
$ws.Range("J11").Value = 0.1504156329559603
$ws.Range("L11").Value = 0.4733220857560383
$ws.Range("N11").Value = 1.392269010674802
$ws.Range("B12").Value = 2.265274486275587
$ws.Range("C12").Value = 0.255906420568607
$ws.Range("D12").Value = 0.1103740471244947
$ws.Range("F12").Value = 1.935416907205862
$ws.Range("G12").Value = 1.331379412021164
$ws.Range("H12").Value = 1.199559386188184
$ws.Range("I12").Value = 1.071125076760481
$ws.Range("J12").Value = 0.1501465157177986
$ws.Range("L12").Value = 0.4760960009393926
$ws.Range("N12").Value = 1.388140145959753
$ws.Range("B13").Value = 2.257898313610724
$ws.Range("C13").Value = 0.2545897976807794
$ws.Range("D13").Value = 0.1103472803775176
$ws.Range("F13").Value = 1.934755890168091
$ws.Range("G13").Value = 1.330553206942568
$ws.Range("H13").Value = 1.199476467690658
$ws.Range("I13").Value = 1.07124751770737
$ws.Range("J13").Value = 0.150204174427861
$ws.Range("L13").Value = 0.4754973756605949
$ws.Range("N13").Value = 1.389025556640288
$ws.Range("B14").Value = 2.233857528318879
$ws.Range("C14").Value = 0.2502949629755733
$ws.Range("D14").Value = 0.1102604354242516
$ws.Range("F14").Value = 1.932628562931427
$ws.Range("G14").Value = 1.327881770450205
$ws.Range("H14").Value = 1.199222034443864
$ws.Range("I14").Value = 1.071662879205526
$ws.Range("J14").Value = 0.150393357092188
$ws.Range("L14").Value = 0.4735497506470381
$ws.Range("N14").Value = 1.391927607867991
$ws.Range("B15").Value = 2.219143029700604
$ws.Range("C15").Value = 0.2476634533449271
$ws.Range("D15").Value = 0.1102075844578465
$ws.Range("F15").Value = 1.931347325527796
$ws.Range("G15").Value = 1.326263101054849
$ws.Range("H15").Value = 1.199078465384844
$ws.Range("I15").Value = 1.071929632865547
$ws.Range("J15").Value = 0.1505101170105156
$ws.Range("L15").Value = 0.4723603279601747
$ws.Range("N15").Value = 1.39371636490894
$ws.Range("B16").Value = 2.134992380286747
$ws.Range("C16").Value = 0.2325711163469748
$ws.Range("D16").Value = 0.10991000787876
$ws.Range("F16").Value = 1.924339153698853
$ws.Range("G16").Value = 1.31725749448654
$ws.Range("H16").Value = 1.198443837069732
$ws.Range("I16").Value = 1.073647145897837
$ws.Range("J16").Value = 0.151192717217743
$ws.Range("L16").Value = 0.4655987870922331
$ws.Range("N16").Value = 1.404138095007635
$ws.Range("B17").Value = 2.083523301268428
$ws.Range("C17").Value = 0.2233015890974741
$ws.Range("D17").Value = 0.1097321942823228
$ws.Range("F17").Value = 1.920337718596485
$ws.Range("G17").Value = 1.31197365076207
$ws.Range("H17").Value = 1.198222267565399
$ws.Range("I17").Value = 1.074869117484553
$ws.Range("J17").Value = 0.1516235211941499
$ws.Range("L17").Value = 0.4614995819224674
$ws.Range("N17").Value = 1.410684076916642
$ws.Range("B18").Value = 2.053976854721952
$ws.Range("C18").Value = 0.2179660228822513
$ws.Range("D18").Value = 0.1096316700152045
$ws.Range("F18").Value = 1.9181455571388
$ws.Range("G18").Value = 1.309022852385056
$ws.Range("H18").Value = 1.198156437279607
$ws.Range("I18").Value = 1.075633751287086
$ws.Range("J18").Value = 0.1518757400218007
$ws.Range("L18").Value = 0.459159792826128
$ws.Range("N18").Value = 1.414505198789612
$ws.Range("B19").Value = 2.0439827974281
$ws.Range("C19").Value = 0.2161588097611116
$ws.Range("D19").Value = 0.1095979353993144
$ws.Range("F19").Value = 1.917422101754923
$ws.Range("G19").Value = 1.308038917601095
$ws.Range("H19").Value = 1.198144726337034
$ws.Range("I19").Value = 1.075903248346343
$ws.Range("J19").Value = 0.1519618984422682
$ws.Range("L19").Value = 0.4583706696489998
$ws.Range("N19").Value = 1.41580858827572
$ws.Range("B20").Value = 2.088996360632393
$ws.Range("C20").Value = 0.2242887582701485
$ws.Range("D20").Value = 0.1097509419556744
$ws.Range("F20").Value = 1.920752357751311
$ws.Range("G20").Value = 1.312526979507965
$ws.Range("H20").Value = 1.198239476677287
$ws.Range("I20").Value = 1.074732640026326
$ws.Range("J20").Value = 0.1515772027738205
$ws.Range("L20").Value = 0.4619340908385112
$ws.Range("N20").Value = 1.409981444737998
$ws.Range("B21").Value = 2.240916562657048
$ws.Range("C21").Value = 0.2515566267014151
$ws.Range("D21").Value = 0.1102858722284452
$ws.Range("F21").Value = 1.933248864348755
$ws.Range("G21").Value = 1.328662754310159
$ws.Range("H21").Value = 1.19929420857838
$ws.Range("I21").Value = 1.071538306918463
$ws.Range("J21").Value = 0.1503376061684527
$ws.Range("L21").Value = 0.4741210749128157
$ws.Range("N21").Value = 1.391072877737347
$ws.Range("B22").Value = 2.340706465011806
$ws.Range("C22").Value = 0.2693420159546633
$ws.Range("D22").Value = 0.1106509000963243
$ws.Range("F22").Value = 1.942392032531501
$ws.Range("G22").Value = 1.339998504856538
$ws.Range("H22").Value = 1.200532996193061
$ws.Range("I22").Value = 1.070002411595262
$ws.Range("J22").Value = 0.1495668565216484
$ws.Range("L22").Value = 0.4822451393847444
$ws.Range("N22").Value = 1.37921475866445
$ws.Range("B23").Value = 2.2874012427244
$ws.Range("C23").Value = 0.2598529036425248
$ws.Range("D23").Value = 0.1104546745378485
$ws.Range("F23").Value = 1.937422744435963
$ws.Range("G23").Value = 1.333875936814536
$ws.Range("H23").Value = 1.199821516063366
$ws.Range("I23").Value = 1.07077158300531
$ws.Range("J23").Value = 0.1499746183266488
$ws.Range("L23").Value = 0.4778946487373332
$ws.Range("N23").Value = 1.385497903373828
$ws.Range("B24").Value = 2.08652185392657
$ws.Range("C24").Value = 0.2238424789972555
$ws.Range("D24").Value = 0.1097424608267232
$ws.Range("F24").Value = 1.920564562003563
$ws.Range("G24").Value = 1.312276548760281
$ws.Range("H24").Value = 1.198231504703585
$ws.Range("I24").Value = 1.074794148038357
$ws.Range("J24").Value = 0.1515981291771826
$ws.Range("L24").Value = 0.4617375967130783
$ws.Range("N24").Value = 1.410298924856384
$ws.Range("B25").Value = 1.872032095094994
$ws.Range("C25").Value = 0.1848473144289358
$ws.Range("D25").Value = 0.1090411983863433
$ws.Range("F25").Value = 1.906559866734554
$ws.Range("G25").Value = 1.292352872175528
$ws.Range("H25").Value = 1.198871398348487
$ws.Range("I25").Value = 1.081494472898569
$ws.Range("J25").Value = 0.1535188279913378
$ws.Range("L25").Value = 0.4449968635679511
$ws.Range("N25").Value = 1.439193337300651
